# Azov city 17-18 (added)
# Fills in previously-blank rows 35/36 with "Азов" data for 2018/2017,
# copies the placeholder ("n/a") formatting from the row above (Сланцы 2014,
# row 34) onto the matching columns, stamps a lone formatted-but-empty cell
# in the new separator row 37, clears out the two now-redundant blank filler
# rows (38/39), and moves the visible selection to match the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the "no data" placeholder formatting pattern from row 34 onto
#        rows 35/36 for the columns that will hold the "???" marker.
#        (Columns G,H,I,J,K,M,N,P,Q,R carry that look in row 34; L,O,S keep
#        the plain numeric style already on rows 35/36.)
$phCols = "G,H,I,J,K,M,N,P,Q,R"
foreach ($col in $phCols.Split(",")) {
    $ws.Range($col + "34").Copy()
    $ws.Range($col + "35").PasteSpecial(-4122)
    $ws.Range($col + "36").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- 2. Row 35: Азов, 2018
$ws.Range("A35").Value = "Азов"
$ws.Range("B35").Value = 2018
$ws.Range("C35").Value = 80.721
$ws.Range("D35").Value = 17.967
$ws.Range("E35").Value = 296
$ws.Range("F35").Value = 29249.1
$ws.Range("G35").Value = "???"
$ws.Range("H35").Value = "???"
$ws.Range("I35").Value = "???"
$ws.Range("J35").Value = "???"
$ws.Range("K35").Value = "???"
$ws.Range("L35").Value = 3672
$ws.Range("M35").Value = "???"
$ws.Range("N35").Value = "???"
$ws.Range("O35").Formula = "= 28823135 / 1000"
$ws.Range("P35").Value = "???"
$ws.Range("Q35").Value = 60
$ws.Range("R35").Value = "???"
$ws.Range("S35").Value = 12216.4
$ws.Range("T35").Value = 413.2
$ws.Range("U35").Value = -51

# --- 3. Row 36: Азов, 2017
$ws.Range("A36").Value = "Азов"
$ws.Range("B36").Value = 2017
$ws.Range("C36").Value = 81.355
$ws.Range("D36").Value = 17.455
$ws.Range("E36").Value = 290
$ws.Range("F36").Value = 27044.2
$ws.Range("G36").Value = "???"
$ws.Range("H36").Value = "???"
$ws.Range("I36").Value = "???"
$ws.Range("J36").Value = "???"
$ws.Range("K36").Value = "???"
$ws.Range("L36").Value = 3050.8
$ws.Range("M36").Value = "???"
$ws.Range("N36").Value = "???"
$ws.Range("O36").Formula = "= 29006627.2/1000"
$ws.Range("P36").Value = "???"
$ws.Range("Q36").Value = 66.543
$ws.Range("R36").Value = "???"
$ws.Range("S36").Value = 11765.9
$ws.Range("T36").Value = 399.8
$ws.Range("U36").Value = -206

# --- 4. Row 37 ends up with a single stray formatted-but-empty cell (B37)
#        carrying the same bold/centered look used by columns A/B elsewhere.
$ws.Range("B35").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 5. The two blank filler rows right below (38/39) are removed outright.
$ws.Rows("38:39").Clear()

# --- 6. Restore the on-screen selection to match the saved view.
$ws.Range("E40").Select()
